$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ball geometry" factor column (column F, originally holding the
# constant 0.5 placeholder for every sport) was dropped from the dataset.
# Deleting the whole column shifts G:P left into F:O and Excel automatically
# drops the now-unused "Ball geometry" shared string.
$ws.Columns("F:F").Delete()

# The last two columns (now N = "Number of rules about movement/Number of
# rules that prevent movement" and O = "UAS") were recomputed after removing
# the factor, so refresh their values per sport.
$ws.Range("N2").Value = 0.2352
$ws.Range("O2").Value = 0.33

$ws.Range("N3").Value = 0.2352
$ws.Range("O3").Value = 0.28

$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.19

$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.23

$ws.Range("N6").Value = 0.5882
$ws.Range("O6").Value = 0.28

$ws.Range("N7").Value = 0.7529
$ws.Range("O7").Value = 0.16

$ws.Range("N8").Value = 0.9411
$ws.Range("O8").Value = 0.06

$ws.Range("N9").Value = 0.7529
$ws.Range("O9").Value = 0.18

$ws.Range("N10").Value = 0.4235
$ws.Range("O10").Value = 0.37

$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0.08

$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0.05

$ws.Range("N13").Value = 0.2352
$ws.Range("O13").Value = 0.3

# Keep the active selection consistent with the edited sheet (matches the
# author re-selecting a cell while reviewing the updated table).
[void]$ws.Range("P10").Select()
